# v1.10: Base histórica id_origen/id_operacion; versión en sidebar; regla versionado al desplegar
$wb = $excel.ActiveWorkbook

# --- Sheet "Log": append bitacora entry as row 41 ---
$log = $wb.Worksheets.Item("Log")
$log.Range("A41").Value = '27/02/2025'
$log.Range("B41").Value = '20:20'
$log.Range("C41").Value = 'Base histórica Excel y versión en sidebar'
$log.Range("D41").Value = 'Export Base histórica: columnas id_origen e id_operacion. Versión de la app visible abajo en el sidebar (APP_VERSION). Regla de bitácora: al indicar desplegar, incrementar versión, actualizar bitácora y desplegar.'
$log.Range("E41").Value = 'Diagnostico'

# --- Sheet "Versiones": append version entry as row 12 ---
$ver = $wb.Worksheets.Item("Versiones")
# Leading apostrophe keeps "1.10" stored as text instead of being parsed as the number 1.1
$ver.Range("A12").Value = "'1.10"
$ver.Range("B12").Value = '27/02/2025'
$ver.Range("C12").Value = 'Base histórica Excel: id_origen e id_operacion; versión en sidebar; regla de versionado al desplegar'
